$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold "Class ID"
$ws.Columns("A:A").Insert()

# Header row
$ws.Range("A1").Value = "Class ID"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "Time"

# Data row 2 (merged with the former row 3 content, with updated date/time).
# Numeric-/date-looking values are forced to text (matching the source
# data's original literal-text formatting) by flipping the number format
# to Text before the assignment, then the cell style is restored to
# Normal afterwards so no stray formatting is left behind.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12345"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1446896"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Asif Newaz"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-01-28"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "08:36:56"
$ws.Range("E2").Style = "Normal"

# Remove the now-unused third row
$ws.Rows("3:3").Delete()
